$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: extend the numbered header sequence (1..10 already in O5:X5) with 11..24 in Y5:AL5,
# copying the formatting used by the existing sequence (style of X5).
$ws.Range("X5").Copy()
$ws.Range("Y5:AL5").PasteSpecial(-4122)
$seq = 11
for ($col = 25; $col -le 38; $col++) {
    $ws.Cells.Item(5, $col).Value = $seq
    $seq++
}

# Rows 6-13: extend the zero-filled columns (already Y:Z = 0) out to column AL,
# copying the formatting used by the existing zero cells (style of Z<row>).
for ($row = 6; $row -le 13; $row++) {
    $ws.Cells.Item($row, 26).Copy()
    $ws.Range($ws.Cells.Item($row, 27), $ws.Cells.Item($row, 38)).PasteSpecial(-4122)
    for ($col = 27; $col -le 38; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

$excel.CutCopyMode = 0

# Update the visible selection to match the new extended range, as Excel does
# when the selection anchor cell's used-range grows.
[void]$ws.Range("O6:AL13").Select()
